# complete link map for phaseII EMTF done, SV generated
#
# Changes:
#  1. The refclk rate for the J11 header link (cells C37:C52 on the
#     "vu13p_gty_refclk" sheet) is doubled: 160.32 -> 320.64.
#  2. Selection / active-sheet bookkeeping left behind by the edit:
#       - "vu13p_gty_refclk" becomes the active (selected) sheet/tab,
#         with C37:C52 highlighted and C37 as the active cell.
#       - "vu13p_gty_inversions" is no longer the active tab (its
#         selection remains anchored on C10).
#       - "Instructions" keeps its selection anchored on B5.

$wb = $excel.ActiveWorkbook

$wsRefclk     = $wb.Worksheets.Item("vu13p_gty_refclk")
$wsInversions = $wb.Worksheets.Item("vu13p_gty_inversions")
$wsInstr      = $wb.Worksheets.Item("Instructions")

# 1. Update the rate values (shared string "160.32" -> "320.64") for the
#    16 rows of the J11 link (C37:C52).
$wsRefclk.Range("C37:C52").Value = "320.64"

# 2. Reproduce the selection state / active tab recorded in the saved
#    workbook: vu13p_gty_refclk is now the active sheet with C37:C52
#    selected (active cell C37).
$wsInversions.Range("C10").Select()
$wsInstr.Range("B5").Select()

$wsRefclk.Activate()
$wsRefclk.Range("C37:C52").Select()
